$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date (column C) for rows 2-9 from 2023-10-22 (45221) to 2023-10-25 (45224)
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45224
}
